# Apply crypto price/volume updates per commit: "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (column D) cells that will hold purely-numeric-looking text to be
# stored as TEXT (matching the original inlineStr string values), not auto-converted numbers.
$ws.Range("D4:D12").NumberFormat = "@"
$ws.Range("D14:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D21:D51").NumberFormat = "@"

# Set the new cell values exactly as they appear in the updated sheet
$ws.Range("D2").Value = '30.598.15'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '1.920.72'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '241.07'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").Value = '0.9987'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '0.4792'
$ws.Range("E7").Value = '  -1.81%  '
$ws.Range("D8").Value = '0.2861'
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("D9").Value = '0.06764'
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").Value = '104.21'
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").Value = '0.07780'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '1.917.90'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '5.254'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").Value = '0.6771'
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").Value = '294.11'
$ws.Range("E16").Value = '  +7.60%  '
$ws.Range("D17").Value = '30.600.31'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = '0.9994'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '0.000007540'
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("D21").Value = '5.471'
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").Value = '0.4681'
$ws.Range("E22").Value = '  -4.77%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '0.9985'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '6.364'
$ws.Range("E24").Value = '  -2.62%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.470'
$ws.Range("E25").Value = '  -3.70%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '167.43'
$ws.Range("E26").Value = '  +0.91%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '19.64'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.114'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.388'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.1001'
$ws.Range("E30").Value = '  -3.65%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.588'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.520'
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.296'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.04769'
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7329'
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.119'
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.710'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01925'
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.625'
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.394'
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '75.11'
$ws.Range("E41").Value = '  -6.35%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '1.987'
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8654'
$ws.Range("E43").Value = '  -4.49%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '106.24'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4294'
$ws.Range("E45").Value = '  -3.22%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9985'
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.472'
$ws.Range("E47").Value = '  -3.69%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '976.15'
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1209'
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '34.92'
$ws.Range("E50").Value = '  -3.46%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '8.929'
$ws.Range("E51").Value = '  -3.39%  '

# Restore default ("Normal") style on the text-forced Price cells so no stray number format lingers
$ws.Range("D4:D12").Style = "Normal"
$ws.Range("D14:D16").Style = "Normal"
$ws.Range("D18:D19").Style = "Normal"
$ws.Range("D21:D51").Style = "Normal"

